$d = $word.ActiveDocument
$t = $d.Tables(1)

# --- Step A: update the first three summary rows to "0M" ---
$t.Rows(1).Cells(1).Range.Text = "0M"
$t.Rows(2).Cells(1).Range.Text = "0M"
$t.Rows(3).Cells(1).Range.Text = "0M"

# --- Step B: collapse the three trailing tab-separated detail rows down
#     to their single leading value (original row indices 34, 35, 36 —
#     these are unaffected by the row insertions done in Step C below,
#     since those insertions happen earlier in the table, at row 4) ---
$t.Rows(34).Cells(1).Range.Text = "100"
$t.Rows(35).Cells(1).Range.Text = "0"
$t.Rows(36).Cells(1).Range.Text = "51"

# --- Step C: insert 10 new single-value rows right after row 3 (i.e.
#     before what is currently row 4), in order ---
$values = @("36", "0.00003", "0.00005", "0.00004", "0.00000", "0.00004", "0.00004", "0.00004", "0.00134", "100.0")
$insertBefore = $t.Rows(4)
for ($i = $values.Length - 1; $i -ge 0; $i--) {
    $newRow = $t.Rows.Add($insertBefore)
    $newRow.Cells(1).Range.Text = $values[$i]
    $insertBefore = $newRow
}
